$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 7, shifting rows 7:432 down to 8:433
$ws.Rows("7:7").Insert()

# Populate the new row 7 with the "enable_year" configuration entry
$ws.Range("A7").Value = "CHE"
$ws.Range("B7").Value = "conv_chp_WASTE"
$ws.Range("C7").Value = "enable_year"
$ws.Range("D7").Value = "configuration"
$ws.Range("G7").Value = 1990

# Re-apply the AutoFilter over the new, larger range (A5:L853)
$ws.AutoFilterMode = $false
foreach ($n in @($wb.Names)) {
    $n.Delete()
}
$ws.Range("A5:L853").AutoFilter()
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$A`$5:`$L`$853")
$filterName.Visible = $false

# Update the active selection to H7
$ws.Range("H7").Select()

Write-Host "Done"
